$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69: Steeling the Knife, Steeling the Mind | Grade 1 Mind Dissolvent
$ws.Range("H69").Value2 = 4014.2
$ws.Range("J69").Value2 = 4081.6667
$ws.Range("L69").Value2 = 12245.0001
$ws.Range("N69").Value2 = -13993.0001

# Row 72: Surgical Substitution (L) | Grade 1 Mind Dissolvent
$ws.Range("H72").Value2 = 4014.2
$ws.Range("J72").Value2 = 4081.6667
$ws.Range("L72").Value2 = 36735.0003
$ws.Range("N72").Value2 = -45471.0003

# Row 115: 5-bell Energy | Competent Craftsman's Syrup
$ws.Range("H115").Value2 = 1910.909
$ws.Range("I115").Value2 = 673.3333
$ws.Range("K115").Value2 = 2019.9999
$ws.Range("M115").Value2 = -452.9999

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value2 = 4547961
$ws.Range("I132").Value2 = 4764276
$ws.Range("K132").Value2 = 14292828
$ws.Range("M132").Value2 = -14290298

# Row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws.Range("H133").Value2 = 27500
$ws.Range("J133").Value2 = 27500
$ws.Range("L133").Value2 = 27500
$ws.Range("N133").Value2 = -37620

# Row 134: Binding Spells | Crocodileskin Index
$ws.Range("H134").Value2 = 0
$ws.Range("J134").Value2 = 0
$ws.Range("L134").Value2 = 0
$ws.Range("N134").ClearContents() | Out-Null

# Row 136: I Like Big Brush and I Cannot Lie | Dark Mahogany Round Brush
$ws.Range("H136").Value2 = 0
$ws.Range("J136").Value2 = 0
$ws.Range("L136").Value2 = 0
$ws.Range("N136").ClearContents() | Out-Null

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value2 = 2639.889
$ws.Range("I137").Value2 = 3098.1785
$ws.Range("J137").Value2 = 1885.0588
$ws.Range("K137").Value2 = 9294.5355
$ws.Range("L137").Value2 = 5655.1764
$ws.Range("M137").Value2 = -6744.5355
$ws.Range("N137").Value2 = -10755.1764

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value2 = 4200.6313
$ws.Range("I138").Value2 = 1367.125
$ws.Range("J138").Value2 = 5508.404
$ws.Range("K138").Value2 = 4101.375
$ws.Range("L138").Value2 = 16525.212
$ws.Range("M138").Value2 = 1038.625
$ws.Range("N138").Value2 = -26805.212

# Row 139: Something Salty and Ceremonial | Gomphotherium Codex
$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").Value2 = 0
$ws.Range("N139").ClearContents() | Out-Null

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value2 = 268698.03
$ws.Range("I141").Value2 = 1058.8286
$ws.Range("K141").Value2 = 3176.4858
$ws.Range("M141").Value2 = 2003.5142

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value2 = 2683.5833
$ws.Range("I61").Value2 = 900.45
$ws.Range("K61").Value2 = 900.45
$ws.Range("M61").Value2 = -688.45

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value2 = 2966.5264
$ws.Range("I122").Value2 = 2345.8333
$ws.Range("J122").Value2 = 4030.5715
$ws.Range("K122").Value2 = 7037.499899999999
$ws.Range("L122").Value2 = 12091.7145
$ws.Range("M122").Value2 = -4587.499899999999
$ws.Range("N122").Value2 = -16991.7145

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value2 = 2683.5833
$ws.Range("I136").Value2 = 900.45
$ws.Range("K136").Value2 = 2701.35
$ws.Range("M136").Value2 = -151.3500000000004

$ws = $wb.Worksheets.Item("BSM")
# Row 19: Twice as Slice | Spiked Bronze Labrys
$ws.Range("H19").Value2 = 9999
$ws.Range("I19").Value2 = 9998
$ws.Range("K19").Value2 = 9998
$ws.Range("M19").Value2 = -9825

# Row 82: Spirituality Inspector | Titanium Lump Hammer
$ws.Range("H82").Value2 = 23409.5
$ws.Range("I82").Value2 = 2381
$ws.Range("J82").Value2 = 38429.855
$ws.Range("K82").Value2 = 2381
$ws.Range("L82").Value2 = 38429.855
$ws.Range("M82").Value2 = -1998
$ws.Range("N82").Value2 = -39195.855

# Row 85: The Clamor for Hammers (L) | Titanium Lump Hammer
$ws.Range("H85").Value2 = 23409.5
$ws.Range("I85").Value2 = 2381
$ws.Range("J85").Value2 = 38429.855
$ws.Range("K85").Value2 = 2381
$ws.Range("L85").Value2 = 38429.855
$ws.Range("M85").Value2 = -1055
$ws.Range("N85").Value2 = -41081.855

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value2 = 2238.5715
$ws.Range("I99").Value2 = 681
$ws.Range("J99").Value2 = 5042.2
$ws.Range("K99").Value2 = 681
$ws.Range("L99").Value2 = 5042.2
$ws.Range("M99").Value2 = 817
$ws.Range("N99").Value2 = -8038.2

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value2 = 2704.4583
$ws.Range("I134").Value2 = 2007.0555
$ws.Range("J134").Value2 = 4796.6665
$ws.Range("K134").Value2 = 6021.166499999999
$ws.Range("L134").Value2 = 14389.9995
$ws.Range("M134").Value2 = -3486.166499999999
$ws.Range("N134").Value2 = -19459.9995

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value2 = 12823373
$ws.Range("I58").Value2 = 1621.4546
$ws.Range("J58").Value2 = 29416228
$ws.Range("K58").Value2 = 1621.4546
$ws.Range("L58").Value2 = 29416228
$ws.Range("M58").Value2 = -1418.4546
$ws.Range("N58").Value2 = -29416634

# Row 68: Do You Even String Bow | Holy Cedar Composite Bow
$ws.Range("H68").Value2 = 32312.857
$ws.Range("J68").Value2 = 32312.857
$ws.Range("L68").Value2 = 32312.857
$ws.Range("N68").Value2 = -33810.857

# Row 71: Win One Bow, Get Three Free (L) | Holy Cedar Composite Bow
$ws.Range("H71").Value2 = 32312.857
$ws.Range("J71").Value2 = 32312.857
$ws.Range("L71").Value2 = 96938.571
$ws.Range("N71").Value2 = -104426.571

# Row 74: License to Heal | Dark Chestnut Rod
$ws.Range("H74").Value2 = 20387.111
$ws.Range("J74").Value2 = 20387.111
$ws.Range("L74").Value2 = 20387.111
$ws.Range("N74").Value2 = -22135.111

# Row 77: Purified Polyrhythm (L) | Dark Chestnut Rod
$ws.Range("H77").Value2 = 20387.111
$ws.Range("J77").Value2 = 20387.111
$ws.Range("L77").Value2 = 61161.333
$ws.Range("N77").Value2 = -69897.333

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value2 = 12823373
$ws.Range("I136").Value2 = 1621.4546
$ws.Range("J136").Value2 = 29416228
$ws.Range("K136").Value2 = 4864.3638
$ws.Range("L136").Value2 = 88248684
$ws.Range("M136").Value2 = -2314.3638
$ws.Range("N136").Value2 = -88253784

$ws = $wb.Worksheets.Item("CUL")
# Row 94: All You Can Stomach | Baklava
$ws.Range("H94").Value2 = 3708.4211
$ws.Range("I94").Value2 = 1475
$ws.Range("J94").Value2 = 3971.1765
$ws.Range("K94").Value2 = 4425
$ws.Range("L94").Value2 = 11913.5295
$ws.Range("M94").Value2 = -3749
$ws.Range("N94").Value2 = -13265.5295

# Row 117: A Good Omen | Peppered Popotoes
$ws.Range("H117").Value2 = 1417.2
$ws.Range("J117").Value2 = 1896
$ws.Range("L117").Value2 = 5688
$ws.Range("N117").Value2 = -12572

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value2 = 1297.4286
$ws.Range("J131").Value2 = 1151.5385
$ws.Range("L131").Value2 = 3454.6155
$ws.Range("N131").Value2 = -13534.6155

# Row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value2 = 2316.95
$ws.Range("I136").Value2 = 1565.2667
$ws.Range("J136").Value2 = 4572
$ws.Range("K136").Value2 = 4695.800099999999
$ws.Range("L136").Value2 = 13716
$ws.Range("M136").Value2 = 404.1999000000005
$ws.Range("N136").Value2 = -23916

# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value2 = 6949716.5
$ws.Range("I139").Value2 = 9617684
$ws.Range("J139").Value2 = 12999.9
$ws.Range("K139").Value2 = 28853052
$ws.Range("L139").Value2 = 38999.7
$ws.Range("M139").Value2 = -28847912
$ws.Range("N139").Value2 = -49279.7

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value2 = 3077.3809
$ws.Range("I122").Value2 = 2278
$ws.Range("J122").Value2 = 4252.9414
$ws.Range("K122").Value2 = 6834
$ws.Range("L122").Value2 = 12758.8242
$ws.Range("M122").Value2 = -4384
$ws.Range("N122").Value2 = -17658.8242

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value2 = 3356.2432
$ws.Range("I132").Value2 = 3058.238
$ws.Range("K132").Value2 = 9174.714
$ws.Range("M132").Value2 = -6644.714

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value2 = 43481610
$ws.Range("I61").Value2 = 62501820
$ws.Range("J61").Value2 = 6842.857
$ws.Range("K61").Value2 = 62501820
$ws.Range("L61").Value2 = 6842.857
$ws.Range("M61").Value2 = -62501618
$ws.Range("N61").Value2 = -7246.857

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value2 = 43481610
$ws.Range("I113").Value2 = 62501820
$ws.Range("J113").Value2 = 6842.857
$ws.Range("K113").Value2 = 62501820
$ws.Range("L113").Value2 = 6842.857
$ws.Range("M113").Value2 = -62499650
$ws.Range("N113").Value2 = -11182.857

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value2 = 1510.8667
$ws.Range("I136").Value2 = 1052.8529
$ws.Range("J136").Value2 = 2926.5454
$ws.Range("K136").Value2 = 3158.5587
$ws.Range("L136").Value2 = 8779.636200000001
$ws.Range("M136").Value2 = -608.5587000000005
$ws.Range("N136").Value2 = -13879.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 43: Walk Softly and Carry a Big Halberd | Velveteen Dress Shoes
$ws.Range("H43").Value2 = 5000
$ws.Range("I43").Value2 = 0
$ws.Range("K43").Value2 = 0
$ws.Range("M43").ClearContents() | Out-Null
